# orders.xlsx: repurpose the "Orders" sheet for a new, simpler header layout
# and drop the sample/test data rows that used to follow it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (A1:J1) replacing the old orderId..timestamp (A1:O1) headers
$ws.Range("A1").Value = "S.No"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Item"
$ws.Range("E1").Value = "Quantity"
$ws.Range("F1").Value = "Amount"
$ws.Range("G1").Value = "Ph no"
$ws.Range("H1").Value = "Tracking ID"
$ws.Range("I1").Value = "Order Status"
$ws.Range("J1").Value = "Timestamp"

# The old sheet had headers/data out to column O; clear what's left beyond
# the new 10-column (A:J) layout.
$ws.Range("K1:O1").Clear()

# Drop the old sample data rows (2-4) entirely, shifting rows up so only
# the header row remains.
$ws.Rows("2:4").Delete()
